$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 15.93194
$ws.Cells.Item(2, 8).Value = 47.79582
$ws.Cells.Item(2, 9).Value = 0.9552847657129105
$ws.Cells.Item(2, 10).Value = 0.9552847657129107
$ws.Cells.Item(2, 13).Value = 31.61061466666667
$ws.Cells.Item(2, 14).Value = 94.831844
$ws.Cells.Item(2, 15).Value = 0.8860472269592234
$ws.Cells.Item(2, 16).Value = 0.8860472269592234
$ws.Cells.Item(2, 17).Value = 503.6184162324533
$ws.Cells.Item(2, 18).Value = 4532.56574609208
$ws.Cells.Item(2, 19).Value = 0.8464274176163158
$ws.Cells.Item(2, 20).Value = 0.8464274176163159

$ws.Cells.Item(3, 7).Value = 15.93194
$ws.Cells.Item(3, 8).Value = 47.79582
$ws.Cells.Item(3, 9).Value = 0.9552847657129105
$ws.Cells.Item(3, 10).Value = 0.9552847657129107
$ws.Cells.Item(3, 15).Value = 0.04688826274109129
$ws.Cells.Item(3, 16).Value = 0.04688826274109129
$ws.Cells.Item(3, 17).Value = 26.65071556354667
$ws.Cells.Item(3, 18).Value = 239.85644007192
$ws.Cells.Item(3, 19).Value = 0.04479164308730879
$ws.Cells.Item(3, 20).Value = 0.04479164308730878

$ws.Cells.Item(4, 7).Value = 15.93194
$ws.Cells.Item(4, 8).Value = 47.79582
$ws.Cells.Item(4, 9).Value = 0.9552847657129105
$ws.Cells.Item(4, 10).Value = 0.9552847657129107
$ws.Cells.Item(4, 13).Value = 2.392593
$ws.Cells.Item(4, 14).Value = 7.177778999999999
$ws.Cells.Item(4, 15).Value = 0.06706451029968528
$ws.Cells.Item(4, 16).Value = 0.06706451029968527
$ws.Cells.Item(4, 17).Value = 38.11864812041999
$ws.Cells.Item(4, 18).Value = 343.0678330837799
$ws.Cells.Item(4, 19).Value = 0.06406570500928593
$ws.Cells.Item(4, 20).Value = 0.06406570500928593

$ws.Cells.Item(5, 9).Value = 0.004609931913019111
$ws.Cells.Item(5, 10).Value = 0.004609931913019112
$ws.Cells.Item(5, 13).Value = 31.61061466666667
$ws.Cells.Item(5, 14).Value = 94.831844
$ws.Cells.Item(5, 15).Value = 0.8860472269592234
$ws.Cells.Item(5, 16).Value = 0.8860472269592234
$ws.Cells.Item(5, 17).Value = 2.430318887417333
$ws.Cells.Item(5, 18).Value = 21.872869986756
$ws.Cells.Item(5, 19).Value = 0.004084617388001412
$ws.Cells.Item(5, 20).Value = 0.004084617388001412

$ws.Cells.Item(6, 9).Value = 0.004609931913019111
$ws.Cells.Item(6, 10).Value = 0.004609931913019112
$ws.Cells.Item(6, 15).Value = 0.04688826274109129
$ws.Cells.Item(6, 16).Value = 0.04688826274109129
$ws.Cells.Item(6, 19).Value = 0.0002161516987561817
$ws.Cells.Item(6, 20).Value = 0.0002161516987561817

$ws.Cells.Item(7, 9).Value = 0.004609931913019111
$ws.Cells.Item(7, 10).Value = 0.004609931913019112
$ws.Cells.Item(7, 13).Value = 2.392593
$ws.Cells.Item(7, 14).Value = 7.177778999999999
$ws.Cells.Item(7, 15).Value = 0.06706451029968528
$ws.Cells.Item(7, 16).Value = 0.06706451029968527
$ws.Cells.Item(7, 17).Value = 0.183949727619
$ws.Cells.Item(7, 18).Value = 1.655547548571
$ws.Cells.Item(7, 19).Value = 0.0003091628262615181
$ws.Cells.Item(7, 20).Value = 0.0003091628262615181

$ws.Cells.Item(8, 7).Value = 0.6688636666666667
$ws.Cells.Item(8, 8).Value = 2.006591
$ws.Cells.Item(8, 9).Value = 0.04010530237407027
$ws.Cells.Item(8, 10).Value = 0.04010530237407027
$ws.Cells.Item(8, 13).Value = 31.61061466666667
$ws.Cells.Item(8, 14).Value = 94.831844
$ws.Cells.Item(8, 15).Value = 0.8860472269592234
$ws.Cells.Item(8, 16).Value = 0.8860472269592234
$ws.Cells.Item(8, 17).Value = 21.14319163153378
$ws.Cells.Item(8, 18).Value = 190.288724683804
$ws.Cells.Item(8, 19).Value = 0.03553519195490612
$ws.Cells.Item(8, 20).Value = 0.03553519195490612

$ws.Cells.Item(9, 7).Value = 0.6688636666666667
$ws.Cells.Item(9, 8).Value = 2.006591
$ws.Cells.Item(9, 9).Value = 0.04010530237407027
$ws.Cells.Item(9, 10).Value = 0.04010530237407027
$ws.Cells.Item(9, 15).Value = 0.04688826274109129
$ws.Cells.Item(9, 16).Value = 0.04688826274109129
$ws.Cells.Item(9, 17).Value = 1.118865331599556
$ws.Cells.Item(9, 18).Value = 10.069787984396
$ws.Cells.Item(9, 19).Value = 0.001880467955026319
$ws.Cells.Item(9, 20).Value = 0.001880467955026319

$ws.Cells.Item(10, 7).Value = 0.6688636666666667
$ws.Cells.Item(10, 8).Value = 2.006591
$ws.Cells.Item(10, 9).Value = 0.04010530237407027
$ws.Cells.Item(10, 10).Value = 0.04010530237407027
$ws.Cells.Item(10, 13).Value = 2.392593
$ws.Cells.Item(10, 14).Value = 7.177778999999999
$ws.Cells.Item(10, 15).Value = 0.06706451029968528
$ws.Cells.Item(10, 16).Value = 0.06706451029968527
$ws.Cells.Item(10, 17).Value = 1.600318526821
$ws.Cells.Item(10, 18).Value = 14.402866741389
$ws.Cells.Item(10, 19).Value = 0.002689642464137828
$ws.Cells.Item(10, 20).Value = 0.002689642464137828

Write-Output "Applied all updates"